# Generate Report for handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de sheets to reflect the new report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-14 03:16:51"
$wsZhCn.Range("G2").Value = "2016-01-14 03:17:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-14 03:17:02"
$wsDeDe.Range("G2").Value = "2016-01-14 03:17:55"
